# Fruta / hortaliza, semanal
# Insert a new weekly record for "Macroferia Regional de Talca" / Kiwi,
# pushing the existing rows 135..188 down to 136..189, then populate the
# newly-opened row 135 with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 135 (shifts 135-188 -> 136-189)
$ws.Rows.Item(135).Insert()

# Fill in the new row 135 with the latest weekly record
$ws.Range("A135").Value = 5
$ws.Range("B135").Value = "Macroferia Regional de Talca"
$ws.Range("C135").Value = "Maule"
$ws.Range("D135").Value = 44466
$ws.Range("E135").Value = 7
$ws.Range("F135").Value = "Fruta"
$ws.Range("G135").Value = 100101
$ws.Range("H135").Value = "Berries"
$ws.Range("I135").Value = 100101007
$ws.Range("J135").Value = "Kiwi"
$ws.Range("K135").Value = "Hayward"
$ws.Range("L135").Value = "Primera"
$ws.Range("M135").Value = 200
$ws.Range("N135").Value = 10000
$ws.Range("O135").Value = 10000
$ws.Range("P135").Value = 10000
$ws.Range("Q135").Value = "$/bandeja 18 kilos"
$ws.Range("R135").Value = "Provincia de Curicó"
$ws.Range("S135").Value = 556
$ws.Range("T135").Value = 18
